# ---------------------------------------------------------------------------
# W2 Assignment / Input Route Data Template.xlsx
# Replace the placeholder single-protest-group template with the real
# 2018-02-02 multi-group march route table:
#   - insert a "신고 인원(명)" (reported headcount) column between the group
#     name and the time column
#   - add a dated sub-title row ("2018.02.02")
#   - populate six march groups with their head counts / times / routes
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe the old layout (content + formatting) in the working area --------
$ws.Range("B1:J20").Clear()

# --- column widths -----------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.625   # B 시위대명
$ws.Columns.Item(3).ColumnWidth = 12.75    # C 신고 인원(명)
$ws.Columns.Item(4).ColumnWidth = 12.125   # D 시간
$ws.Columns.Item(5).ColumnWidth = 86.25    # E 행진 경로(교차로)
$ws.Columns.Item(6).ColumnWidth = 11.875   # F (spacer)

# --- row heights ---------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 12
$ws.Rows.Item(2).RowHeight = 20.25
$ws.Rows.Item(3).RowHeight = 11.25
$ws.Rows.Item(4).RowHeight = 21
$ws.Rows.Item(6).RowHeight = 29.25
$ws.Rows.Item(7).RowHeight = 18

# ===========================================================================
# Values first (text cells forced via NumberFormat="@" before the write so
# things like "2018.02.02" don't get reinterpreted as a date serial)
# ===========================================================================
$ws.Range("B1:F4").NumberFormat = "@"

$ws.Range("B2").Value = "중부, 종로, 남대문 관내 행진 교차로 Template"
$ws.Range("B4").Value = "2018.02.02"

$ws.Range("B5").Value = "시위대명"
$ws.Range("C5").Value = "신고 인원(명)"
$ws.Range("D5").Value = "시간"
$ws.Range("E5").Value = "행진 경로(교차로)"

$data = @(
    @("새한국", 300, "15:40~17:00", "청계광장→세종→세문관→광화문→동십자→안국→인사동→종로2→종로1→서린→세종→청계광장"),
    @("구명총", 100, "15:30~17:30", "세종↔서린↔종로1↔종로2↔종로3"),
    @("태극기국민평의회", 100, "15:30~17:30", "한은→눈스퀘어→을지1→광교→종로1→서린→세종"),
    @("석방운동본부", 3000, "15:30~19:00", "서울역→남대문→한은→눈스퀘어→을지1→광교→종로1→서린→세종→세문관↔광화문"),
    @("태극기행동본부", 300, "15:20~18:00", "세종→세문관→광화문→동십자→안국→인사동→종로2→종로1→서린→세종"),
    @("태극기국민운동본부", 600, "15:30~18:00", "대한문→환구단→개풍→을지1→눈스퀘어→한은→남대문→태평→대한문")
)

$row = 6
foreach ($entry in $data) {
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $entry[0]

    # number stays numeric: write value first, THEN stamp the text format
    $ws.Range("C$row").Value = $entry[1]
    $ws.Range("C$row").NumberFormat = "@"

    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $entry[2]

    $ws.Range("E$row").NumberFormat = "@"
    $ws.Range("E$row").Value = $entry[3]

    $row = $row + 1
}

# ===========================================================================
# Merges
# ===========================================================================
$ws.Range("B2:E2").Merge()
$ws.Range("B4:E4").Merge()

# ===========================================================================
# Fonts (never touch .Font.Name explicitly -- the sheet's default font is
# already 맑은 고딕, and re-stamping the Name property mints a redundant
# near-duplicate font entry that drops the theme/minor-scheme linkage).
# The bold header/title fonts in this template are all "family 3"; stamp
# that explicitly or a plain Bold+Size mutation mints a stray "family 2"
# duplicate of the same face instead of reusing the canonical one.
# ===========================================================================
$ws.Range("B2:E2").Font.Bold = $true
$ws.Range("B2:E2").Font.Size = 14
$ws.Range("B2:E2").Font.Family = 3

$ws.Range("B3:E3").Font.Bold = $true
$ws.Range("B3:E3").Font.Size = 11
$ws.Range("B3:E3").Font.Family = 3

$ws.Range("B4:E4").Font.Bold = $true
$ws.Range("B4:E4").Font.Size = 12
$ws.Range("B4:E4").Font.Family = 3

$ws.Range("B5:E5").Font.Bold = $true
$ws.Range("B5:E5").Font.Size = 11
$ws.Range("B5:E5").Font.Family = 3

$ws.Range("B6:E11").Font.Bold = $false
$ws.Range("B6:E11").Font.Size = 11

# ===========================================================================
# Alignment
# ===========================================================================
$ws.Range("B1:F4").HorizontalAlignment = -4108  # xlCenter
$ws.Range("B1:F4").VerticalAlignment = -4108

$ws.Range("B4:E4").HorizontalAlignment = -4152  # xlRight (overrides the block above)

$ws.Range("B5:E11").HorizontalAlignment = -4108
$ws.Range("B5:E11").VerticalAlignment = -4108

# ===========================================================================
# Borders
# ===========================================================================
# Row 2 title block: left+top+bottom framing (no inner/right divider)
$b2 = $ws.Range("B2")
$b2.Borders.Item(7).LineStyle = 1
$b2.Borders.Item(7).Weight = 2
$b2.Borders.Item(8).LineStyle = 1
$b2.Borders.Item(8).Weight = 2
$b2.Borders.Item(9).LineStyle = 1
$b2.Borders.Item(9).Weight = 2

$ce2 = $ws.Range("C2:E2")
$ce2.Borders.Item(8).LineStyle = 1
$ce2.Borders.Item(8).Weight = 2
$ce2.Borders.Item(9).LineStyle = 1
$ce2.Borders.Item(9).Weight = 2

# Row 4 date block: bottom framing everywhere, left/right dividers inside
$b4 = $ws.Range("B4")
$b4.Borders.Item(9).LineStyle = 1
$b4.Borders.Item(9).Weight = 2
$b4.Borders.Item(10).LineStyle = 1
$b4.Borders.Item(10).Weight = 2

$cd4 = $ws.Range("C4:D4")
$cd4.Borders.Item(7).LineStyle = 1
$cd4.Borders.Item(7).Weight = 2
$cd4.Borders.Item(9).LineStyle = 1
$cd4.Borders.Item(9).Weight = 2
$cd4.Borders.Item(10).LineStyle = 1
$cd4.Borders.Item(10).Weight = 2

$e4 = $ws.Range("E4")
$e4.Borders.Item(7).LineStyle = 1
$e4.Borders.Item(7).Weight = 2
$e4.Borders.Item(9).LineStyle = 1
$e4.Borders.Item(9).Weight = 2

# Header row + data grid: thin box around every cell
$ws.Range("B5:E11").Borders.LineStyle = 1
$ws.Range("B5:E11").Borders.Weight = 2

# ===========================================================================
# Misc view tweaks
# ===========================================================================
$ws.Range("D14").Select()
